# generator/variableOM.xlsx -- "update generator data to NETL data"
#
# Row 2 used to hold the *labels* of the technology columns (as text,
# re-using shared strings) instead of real numbers. Replace those with the
# actual NETL data values, add a new helper row 5 (ratio derived from
# ES-4, used to compute the ng-cc-new value in E2), and tidy up the
# view/selection the way the author left the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new helper row 5: ratio pulled from ES-4 ---------------------------
$ws.Range("A5").Value2 = "ratio from ES-4"
$ws.Range("B5").Value2 = 1.5
$ws.Range("F5").Formula = "=5.6/1.7"

# --- row 2: replace placeholder text with real NETL numbers -------------
$ws.Range("B2").Value2 = 7.5
$ws.Range("C2").Value2 = 5
$ws.Range("E2").Formula = "=F2*F5"
$ws.Range("F2").Value2 = 4
$ws.Range("G2").Value2 = 5
$ws.Range("H2").Value2 = 2

# --- cosmetic touch-ups matching the author's saved view -----------------
$ws.Columns.Item(1).EntireColumn.AutoFit() | Out-Null
$ws.Range("F6").Select() | Out-Null
